# Generate Report for Handoff
#
# 1. The handoff status text changes everywhere it appears ("Ready for
#    handoff" -> "Handoff transform failed") -- this shows up on the
#    Overview sheet (columns B & C) as well as on each locale sheet's
#    row 2 (column B), since they all shared the same text.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "Handoff transform failed", 1)
}

# 2. Each locale sheet (zh-cn, de-de) had a successful-handoff row (row 2)
#    pointing at a generated .xlf file with a real handoff datetime and a
#    "Include" dependency reason. The transform now failed, so: the
#    handoff-file link/cell is removed, the handoff datetime collapses to
#    the "never happened" sentinel date, and the reason flips to
#    "Ignored".
foreach ($name in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($name)

    # Remove the "Latest Handoff File" cell (and its hyperlink) entirely.
    $ws.Range("C2").ClearContents()

    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"

    # ClearContents() does not drop the hyperlink anchored to C2, and this
    # object model only supports removing hyperlinks at the collection
    # level (Hyperlinks.Delete() removes every link on the sheet). Rebuild
    # the two links that must survive (the source-file link on A2 and the
    # .localization-config link on A3).
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f7731970b1100e44d06028a792caa070a7655bf8/e2e/3c182483-3a20-4aae-a6e6-7d025a8a4222.md", [System.Type]::Missing, [System.Type]::Missing, "3c182483-3a20-4aae-a6e6-7d025a8a4222.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f7731970b1100e44d06028a792caa070a7655bf8/.localization-config", [System.Type]::Missing, [System.Type]::Missing, ".localization-config")
}
